$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting FIRST NAME / LAST NAME / MIDDLE NAME / CATEGORY right
$ws.Columns.Item(1).Insert()

# New column A header
$ws.Range("A1").Value = "CUSTOMER ID"

# New column A data
$ws.Range("A2").Value = "cust 01"
$ws.Range("A3").Value = "cust 02"
$ws.Range("A4").Value = "cust 03"
$ws.Range("A5").Value = "cust 04"
$ws.Range("A6").Value = "cust 05"
$ws.Range("A7").Value = "cust 06"
$ws.Range("A8").Value = "cust 07"

# Match the bestFit column width used by Excel for this content
$ws.Columns.Item(1).ColumnWidth = 11.75

# Match the resulting selection state
$ws.Range("C12").Select()
